$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date strings in column A from DD/MM/YYYY to DD-MM-YYYY format.
#
# Plain assignment of a dash-separated date string (e.g. "01-08-2022") gets
# auto-recognised by Excel as an actual date whenever the leading number
# could plausibly be a month (<= 12), which would turn the cell into a
# numeric/date cell instead of leaving it as literal text. To avoid that,
# such values are entered with a leading apostrophe (forcing text entry)
# and then ClearFormats() removes the resulting quote-prefix cell style so
# the cell ends up with no explicit style applied, matching a plain text
# entry like the unambiguous dates get by default.
$ambiguousRows = @(4, 5, 6, 7, 13, 14, 15, 16)
for ($row = 3; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $newDate = $cell.Value2 -replace '/', '-'
    if ($ambiguousRows -contains $row) {
        $cell.Value = "'" + $newDate
        $cell.ClearFormats()
    } else {
        $cell.Value = $newDate
    }
}

# Rows where attendance counts changed: D and E go 0 -> 1, H goes 1 -> 0
$rowsToFlip = @(4, 5, 10, 11, 12)
foreach ($row in $rowsToFlip) {
    $ws.Cells.Item($row, 4).Value = 1   # D column
    $ws.Cells.Item($row, 5).Value = 1   # E column
    $ws.Cells.Item($row, 8).Value = 0   # H column
}
